# Auto-generated edit script
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap F:V content between row pairs (home/away match order corrections) ---
# Row 5
$ws.Cells.Item(5, 6).Value = 'MC Alger'
$ws.Cells.Item(5, 7).Value = 4
$ws.Cells.Item(5, 8).Value = 'Ben Aknoun'
$ws.Cells.Item(5, 9).Value = 0
$ws.Cells.Item(5, 10).Value = 1.4
$ws.Cells.Item(5, 11).Value = '16/09/2023 03:43'
$ws.Cells.Item(5, 12).Value = 1.3
$ws.Cells.Item(5, 13).Value = '16/09/2023 10:40'
$ws.Cells.Item(5, 14).Value = 4.19
$ws.Cells.Item(5, 15).Value = '16/09/2023 03:43'
$ws.Cells.Item(5, 16).Value = 4.81
$ws.Cells.Item(5, 17).Value = '16/09/2023 16:47'
$ws.Cells.Item(5, 18).Value = 8.529999999999999
$ws.Cells.Item(5, 19).Value = '16/09/2023 03:43'
$ws.Cells.Item(5, 20).Value = 12.64
$ws.Cells.Item(5, 21).Value = '16/09/2023 16:47'
$ws.Cells.Item(5, 22).Value = 'https://www.betexplorer.com/football/algeria/ligue-1/mc-alger-es-ben-aknoun/WjyqCu9h/'

# Row 6
$ws.Cells.Item(6, 6).Value = 'Magra'
$ws.Cells.Item(6, 7).Value = 0
$ws.Cells.Item(6, 8).Value = 'Kabylie'
$ws.Cells.Item(6, 9).Value = 1
$ws.Cells.Item(6, 10).Value = 2.84
$ws.Cells.Item(6, 11).Value = '15/09/2023 13:42'
$ws.Cells.Item(6, 12).Value = 3.8
$ws.Cells.Item(6, 13).Value = '16/09/2023 16:12'
$ws.Cells.Item(6, 14).Value = 2.63
$ws.Cells.Item(6, 15).Value = '15/09/2023 13:42'
$ws.Cells.Item(6, 16).Value = 2.84
$ws.Cells.Item(6, 17).Value = '16/09/2023 15:03'
$ws.Cells.Item(6, 18).Value = 2.72
$ws.Cells.Item(6, 19).Value = '15/09/2023 13:42'
$ws.Cells.Item(6, 20).Value = 2.26
$ws.Cells.Item(6, 21).Value = '16/09/2023 16:12'
$ws.Cells.Item(6, 22).Value = 'https://www.betexplorer.com/football/algeria/ligue-1/magra-kabylie/YFXa8c8H/'

# Row 8
$ws.Cells.Item(8, 6).Value = 'Ben Aknoun'
$ws.Cells.Item(8, 7).Value = 2
$ws.Cells.Item(8, 8).Value = 'ASO Chlef'
$ws.Cells.Item(8, 9).Value = 2
$ws.Cells.Item(8, 10).Value = 2.81
$ws.Cells.Item(8, 11).Value = '22/09/2023 08:13'
$ws.Cells.Item(8, 12).Value = 2.51
$ws.Cells.Item(8, 13).Value = '22/09/2023 16:00'
$ws.Cells.Item(8, 14).Value = 3.04
$ws.Cells.Item(8, 15).Value = '22/09/2023 08:13'
$ws.Cells.Item(8, 16).Value = 3.06
$ws.Cells.Item(8, 17).Value = '22/09/2023 16:34'
$ws.Cells.Item(8, 18).Value = 2.64
$ws.Cells.Item(8, 19).Value = '22/09/2023 08:13'
$ws.Cells.Item(8, 20).Value = 3.01
$ws.Cells.Item(8, 21).Value = '22/09/2023 16:00'
$ws.Cells.Item(8, 22).Value = 'https://www.betexplorer.com/football/algeria/ligue-1/es-ben-aknoun-aso-chlef/WMgbNNKE/'

# Row 9
$ws.Cells.Item(9, 6).Value = 'Constantine'
$ws.Cells.Item(9, 7).Value = 2
$ws.Cells.Item(9, 8).Value = 'MC Alger'
$ws.Cells.Item(9, 9).Value = 1
$ws.Cells.Item(9, 10).Value = 1.98
$ws.Cells.Item(9, 11).Value = '21/09/2023 05:12'
$ws.Cells.Item(9, 12).Value = 3.68
$ws.Cells.Item(9, 13).Value = '22/09/2023 16:51'
$ws.Cells.Item(9, 14).Value = 2.89
$ws.Cells.Item(9, 15).Value = '21/09/2023 05:12'
$ws.Cells.Item(9, 16).Value = 2.92
$ws.Cells.Item(9, 17).Value = '22/09/2023 16:51'
$ws.Cells.Item(9, 18).Value = 3.97
$ws.Cells.Item(9, 19).Value = '21/09/2023 05:12'
$ws.Cells.Item(9, 20).Value = 2.26
$ws.Cells.Item(9, 21).Value = '22/09/2023 16:51'
$ws.Cells.Item(9, 22).Value = 'https://www.betexplorer.com/football/algeria/ligue-1/constantine-mc-alger/Eoq3MszL/'

# Row 16
$ws.Cells.Item(16, 6).Value = 'US Souf'
$ws.Cells.Item(16, 7).Value = 0
$ws.Cells.Item(16, 8).Value = 'Oran'
$ws.Cells.Item(16, 9).Value = 0
$ws.Cells.Item(16, 10).Value = 2.49
$ws.Cells.Item(16, 11).Value = '28/09/2023 19:27'
$ws.Cells.Item(16, 12).Value = 2.14
$ws.Cells.Item(16, 13).Value = '29/09/2023 13:29'
$ws.Cells.Item(16, 14).Value = 2.88
$ws.Cells.Item(16, 15).Value = '28/09/2023 19:27'
$ws.Cells.Item(16, 16).Value = 2.74
$ws.Cells.Item(16, 17).Value = '29/09/2023 14:49'
$ws.Cells.Item(16, 18).Value = 3.18
$ws.Cells.Item(16, 19).Value = '28/09/2023 19:27'
$ws.Cells.Item(16, 20).Value = 4.43
$ws.Cells.Item(16, 21).Value = '29/09/2023 15:47'
$ws.Cells.Item(16, 22).Value = 'https://www.betexplorer.com/football/algeria/ligue-1/us-souf-oran/6qOsFaSf/'

# Row 18
$ws.Cells.Item(18, 6).Value = 'Khenchela'
$ws.Cells.Item(18, 7).Value = 2
$ws.Cells.Item(18, 8).Value = 'Kabylie'
$ws.Cells.Item(18, 9).Value = 1
$ws.Cells.Item(18, 10).Value = 2.63
$ws.Cells.Item(18, 11).Value = '28/09/2023 04:12'
$ws.Cells.Item(18, 12).Value = 2.05
$ws.Cells.Item(18, 13).Value = '29/09/2023 16:41'
$ws.Cells.Item(18, 14).Value = 2.62
$ws.Cells.Item(18, 15).Value = '28/09/2023 04:12'
$ws.Cells.Item(18, 16).Value = 2.75
$ws.Cells.Item(18, 17).Value = '29/09/2023 16:41'
$ws.Cells.Item(18, 18).Value = 3.02
$ws.Cells.Item(18, 19).Value = '28/09/2023 04:12'
$ws.Cells.Item(18, 20).Value = 4.88
$ws.Cells.Item(18, 21).Value = '29/09/2023 16:27'
$ws.Cells.Item(18, 22).Value = 'https://www.betexplorer.com/football/algeria/ligue-1/khenchela-kabylie/pUZYGLcr/'

# Row 23
$ws.Cells.Item(23, 6).Value = 'Oran'
$ws.Cells.Item(23, 7).Value = 0
$ws.Cells.Item(23, 8).Value = 'Magra'
$ws.Cells.Item(23, 9).Value = 0
$ws.Cells.Item(23, 10).Value = 1.98
$ws.Cells.Item(23, 11).Value = '05/10/2023 07:24'
$ws.Cells.Item(23, 12).Value = 1.65
$ws.Cells.Item(23, 13).Value = '06/10/2023 17:55'
$ws.Cells.Item(23, 14).Value = 2.89
$ws.Cells.Item(23, 15).Value = '05/10/2023 07:24'
$ws.Cells.Item(23, 16).Value = 3.28
$ws.Cells.Item(23, 17).Value = '06/10/2023 19:03'
$ws.Cells.Item(23, 18).Value = 3.98
$ws.Cells.Item(23, 19).Value = '05/10/2023 07:24'
$ws.Cells.Item(23, 20).Value = 5.91
$ws.Cells.Item(23, 21).Value = '06/10/2023 17:55'
$ws.Cells.Item(23, 22).Value = 'https://www.betexplorer.com/football/algeria/ligue-1/oran-magra/WrVYZ04K/'

# Row 24
$ws.Cells.Item(24, 6).Value = 'Biskra'
$ws.Cells.Item(24, 7).Value = 0
$ws.Cells.Item(24, 8).Value = 'Paradou'
$ws.Cells.Item(24, 9).Value = 5
$ws.Cells.Item(24, 10).Value = 2.18
$ws.Cells.Item(24, 11).Value = '05/10/2023 07:24'
$ws.Cells.Item(24, 12).Value = 2.03
$ws.Cells.Item(24, 13).Value = '06/10/2023 19:34'
$ws.Cells.Item(24, 14).Value = 2.84
$ws.Cells.Item(24, 15).Value = '05/10/2023 07:24'
$ws.Cells.Item(24, 16).Value = 3
$ws.Cells.Item(24, 17).Value = '06/10/2023 18:05'
$ws.Cells.Item(24, 18).Value = 3.56
$ws.Cells.Item(24, 19).Value = '05/10/2023 07:24'
$ws.Cells.Item(24, 20).Value = 4.33
$ws.Cells.Item(24, 21).Value = '06/10/2023 19:34'
$ws.Cells.Item(24, 22).Value = 'https://www.betexplorer.com/football/algeria/ligue-1/biskra-paradou/hhWUzskE/'

# Row 41
$ws.Cells.Item(41, 6).Value = 'El Bayadh'
$ws.Cells.Item(41, 7).Value = 4
$ws.Cells.Item(41, 8).Value = 'US Souf'
$ws.Cells.Item(41, 9).Value = 0
$ws.Cells.Item(41, 10).Value = 1.45
$ws.Cells.Item(41, 11).Value = '16/11/2023 03:42'
$ws.Cells.Item(41, 12).Value = 1.45
$ws.Cells.Item(41, 13).Value = '17/11/2023 15:18'
$ws.Cells.Item(41, 14).Value = 3.74
$ws.Cells.Item(41, 15).Value = '16/11/2023 03:42'
$ws.Cells.Item(41, 16).Value = 4.01
$ws.Cells.Item(41, 17).Value = '17/11/2023 15:18'
$ws.Cells.Item(41, 18).Value = 7.38
$ws.Cells.Item(41, 19).Value = '16/11/2023 03:42'
$ws.Cells.Item(41, 20).Value = 8.609999999999999
$ws.Cells.Item(41, 21).Value = '17/11/2023 15:18'
$ws.Cells.Item(41, 22).Value = 'https://www.betexplorer.com/football/algeria/ligue-1/el-bayadh-us-souf/UqfmFEv0/'

# Row 42
$ws.Cells.Item(42, 6).Value = 'Constantine'
$ws.Cells.Item(42, 7).Value = 0
$ws.Cells.Item(42, 8).Value = 'Magra'
$ws.Cells.Item(42, 9).Value = 1
$ws.Cells.Item(42, 10).Value = 1.59
$ws.Cells.Item(42, 11).Value = '26/10/2023 04:42'
$ws.Cells.Item(42, 12).Value = 1.31
$ws.Cells.Item(42, 13).Value = '17/11/2023 15:24'
$ws.Cells.Item(42, 14).Value = 3.47
$ws.Cells.Item(42, 15).Value = '26/10/2023 04:42'
$ws.Cells.Item(42, 16).Value = 4.98
$ws.Cells.Item(42, 17).Value = '17/11/2023 15:24'
$ws.Cells.Item(42, 18).Value = 5.42
$ws.Cells.Item(42, 19).Value = '26/10/2023 04:42'
$ws.Cells.Item(42, 20).Value = 10.79
$ws.Cells.Item(42, 21).Value = '17/11/2023 15:24'
$ws.Cells.Item(42, 22).Value = 'https://www.betexplorer.com/football/algeria/ligue-1/constantine-magra/jgeqGfPg/'

# Row 51
$ws.Cells.Item(51, 6).Value = 'US Souf'
$ws.Cells.Item(51, 7).Value = 3
$ws.Cells.Item(51, 8).Value = 'MC Alger'
$ws.Cells.Item(51, 9).Value = 4
$ws.Cells.Item(51, 10).Value = 4.45
$ws.Cells.Item(51, 11).Value = '24/11/2023 03:13'
$ws.Cells.Item(51, 12).Value = 5.84
$ws.Cells.Item(51, 13).Value = '25/11/2023 14:56'
$ws.Cells.Item(51, 14).Value = 2.95
$ws.Cells.Item(51, 15).Value = '24/11/2023 03:13'
$ws.Cells.Item(51, 16).Value = 3.36
$ws.Cells.Item(51, 17).Value = '25/11/2023 14:56'
$ws.Cells.Item(51, 18).Value = 1.89
$ws.Cells.Item(51, 19).Value = '24/11/2023 03:13'
$ws.Cells.Item(51, 20).Value = 1.7
$ws.Cells.Item(51, 21).Value = '25/11/2023 14:56'
$ws.Cells.Item(51, 22).Value = 'https://www.betexplorer.com/football/algeria/ligue-1/us-souf-mc-alger/Uyh5oKM9/'

# Row 52
$ws.Cells.Item(52, 6).Value = 'Khenchela'
$ws.Cells.Item(52, 7).Value = 0
$ws.Cells.Item(52, 8).Value = 'Ben Aknoun'
$ws.Cells.Item(52, 9).Value = 0
$ws.Cells.Item(52, 10).Value = 1.36
$ws.Cells.Item(52, 11).Value = '24/11/2023 03:13'
$ws.Cells.Item(52, 12).Value = 1.32
$ws.Cells.Item(52, 13).Value = '25/11/2023 14:56'
$ws.Cells.Item(52, 14).Value = 4.19
$ws.Cells.Item(52, 15).Value = '24/11/2023 03:13'
$ws.Cells.Item(52, 16).Value = 4.72
$ws.Cells.Item(52, 17).Value = '25/11/2023 14:56'
$ws.Cells.Item(52, 18).Value = 7.76
$ws.Cells.Item(52, 19).Value = '24/11/2023 03:13'
$ws.Cells.Item(52, 20).Value = 11.86
$ws.Cells.Item(52, 21).Value = '25/11/2023 14:56'
$ws.Cells.Item(52, 22).Value = 'https://www.betexplorer.com/football/algeria/ligue-1/khenchela-es-ben-aknoun/bwyflMxi/'

# --- Append new rows 75-81 (copy formatting from row 74, then set values) ---
# Row 75
$ws.Range("A74:V74").Copy()
$ws.Range("A75:V75").PasteSpecial(-4122)
$ws.Cells.Item(75, 1).Value = 74
$ws.Cells.Item(75, 2).Value = 'algeria'
$ws.Cells.Item(75, 3).Value = 'ligue-1'
$ws.Cells.Item(75, 4).Value = '2023-2024'
$ws.Cells.Item(75, 5).Value = 45288.79166666666
$ws.Cells.Item(75, 6).Value = 'CR Belouizdad'
$ws.Cells.Item(75, 7).Value = 2
$ws.Cells.Item(75, 8).Value = 'Constantine'
$ws.Cells.Item(75, 9).Value = 1
$ws.Cells.Item(75, 10).Value = 1.77
$ws.Cells.Item(75, 11).Value = '27/12/2023 07:12'
$ws.Cells.Item(75, 12).Value = 1.57
$ws.Cells.Item(75, 13).Value = '28/12/2023 18:32'
$ws.Cells.Item(75, 14).Value = 3.12
$ws.Cells.Item(75, 15).Value = '27/12/2023 07:12'
$ws.Cells.Item(75, 16).Value = 3.57
$ws.Cells.Item(75, 17).Value = '28/12/2023 18:32'
$ws.Cells.Item(75, 18).Value = 4.87
$ws.Cells.Item(75, 19).Value = '27/12/2023 07:12'
$ws.Cells.Item(75, 20).Value = 7.23
$ws.Cells.Item(75, 21).Value = '28/12/2023 18:32'
$ws.Cells.Item(75, 22).Value = 'https://www.betexplorer.com/football/algeria/ligue-1/cr-belouizdad-constantine/tp3QLYIb/'

# Row 76
$ws.Range("A74:V74").Copy()
$ws.Range("A76:V76").PasteSpecial(-4122)
$ws.Cells.Item(76, 1).Value = 75
$ws.Cells.Item(76, 2).Value = 'algeria'
$ws.Cells.Item(76, 3).Value = 'ligue-1'
$ws.Cells.Item(76, 4).Value = '2023-2024'
$ws.Cells.Item(76, 5).Value = 45289.63541666666
$ws.Cells.Item(76, 6).Value = 'Khenchela'
$ws.Cells.Item(76, 7).Value = 2
$ws.Cells.Item(76, 8).Value = 'ASO Chlef'
$ws.Cells.Item(76, 9).Value = 1
$ws.Cells.Item(76, 10).Value = 2.11
$ws.Cells.Item(76, 11).Value = '21/12/2023 03:42'
$ws.Cells.Item(76, 12).Value = 1.83
$ws.Cells.Item(76, 13).Value = '29/12/2023 14:40'
$ws.Cells.Item(76, 14).Value = 2.8
$ws.Cells.Item(76, 15).Value = '21/12/2023 03:42'
$ws.Cells.Item(76, 16).Value = 3.22
$ws.Cells.Item(76, 17).Value = '29/12/2023 14:40'
$ws.Cells.Item(76, 18).Value = 3.83
$ws.Cells.Item(76, 19).Value = '21/12/2023 03:42'
$ws.Cells.Item(76, 20).Value = 4.96
$ws.Cells.Item(76, 21).Value = '29/12/2023 14:40'
$ws.Cells.Item(76, 22).Value = 'https://www.betexplorer.com/football/algeria/ligue-1/khenchela-aso-chlef/G2V34XmO/'

# Row 77
$ws.Range("A74:V74").Copy()
$ws.Range("A77:V77").PasteSpecial(-4122)
$ws.Cells.Item(77, 1).Value = 76
$ws.Cells.Item(77, 2).Value = 'algeria'
$ws.Cells.Item(77, 3).Value = 'ligue-1'
$ws.Cells.Item(77, 4).Value = '2023-2024'
$ws.Cells.Item(77, 5).Value = 45289.63541666666
$ws.Cells.Item(77, 6).Value = 'Magra'
$ws.Cells.Item(77, 7).Value = 1
$ws.Cells.Item(77, 8).Value = 'US Souf'
$ws.Cells.Item(77, 9).Value = 0
$ws.Cells.Item(77, 10).Value = 1.65
$ws.Cells.Item(77, 11).Value = '28/12/2023 07:12'
$ws.Cells.Item(77, 12).Value = 1.72
$ws.Cells.Item(77, 13).Value = '29/12/2023 15:07'
$ws.Cells.Item(77, 14).Value = 3.28
$ws.Cells.Item(77, 15).Value = '28/12/2023 07:12'
$ws.Cells.Item(77, 16).Value = 3.42
$ws.Cells.Item(77, 17).Value = '29/12/2023 15:07'
$ws.Cells.Item(77, 18).Value = 5.55
$ws.Cells.Item(77, 19).Value = '28/12/2023 07:12'
$ws.Cells.Item(77, 20).Value = 5.41
$ws.Cells.Item(77, 21).Value = '29/12/2023 15:07'
$ws.Cells.Item(77, 22).Value = 'https://www.betexplorer.com/football/algeria/ligue-1/magra-us-souf/6RADOzZu/'

# Row 78
$ws.Range("A74:V74").Copy()
$ws.Range("A78:V78").PasteSpecial(-4122)
$ws.Cells.Item(78, 1).Value = 77
$ws.Cells.Item(78, 2).Value = 'algeria'
$ws.Cells.Item(78, 3).Value = 'ligue-1'
$ws.Cells.Item(78, 4).Value = '2023-2024'
$ws.Cells.Item(78, 5).Value = 45289.63541666666
$ws.Cells.Item(78, 6).Value = 'Paradou'
$ws.Cells.Item(78, 7).Value = 0
$ws.Cells.Item(78, 8).Value = 'Saoura'
$ws.Cells.Item(78, 9).Value = 0
$ws.Cells.Item(78, 10).Value = 2.12
$ws.Cells.Item(78, 11).Value = '21/12/2023 03:42'
$ws.Cells.Item(78, 12).Value = 2.19
$ws.Cells.Item(78, 13).Value = '29/12/2023 15:10'
$ws.Cells.Item(78, 14).Value = 2.8
$ws.Cells.Item(78, 15).Value = '21/12/2023 03:42'
$ws.Cells.Item(78, 16).Value = 2.84
$ws.Cells.Item(78, 17).Value = '29/12/2023 15:10'
$ws.Cells.Item(78, 18).Value = 3.77
$ws.Cells.Item(78, 19).Value = '21/12/2023 03:42'
$ws.Cells.Item(78, 20).Value = 4.02
$ws.Cells.Item(78, 21).Value = '29/12/2023 15:08'
$ws.Cells.Item(78, 22).Value = 'https://www.betexplorer.com/football/algeria/ligue-1/paradou-saoura/fgU73i3U/'

# Row 79
$ws.Range("A74:V74").Copy()
$ws.Range("A79:V79").PasteSpecial(-4122)
$ws.Cells.Item(79, 1).Value = 78
$ws.Cells.Item(79, 2).Value = 'algeria'
$ws.Cells.Item(79, 3).Value = 'ligue-1'
$ws.Cells.Item(79, 4).Value = '2023-2024'
$ws.Cells.Item(79, 5).Value = 45289.66666666666
$ws.Cells.Item(79, 6).Value = 'ES Setif'
$ws.Cells.Item(79, 7).Value = 1
$ws.Cells.Item(79, 8).Value = 'Oran'
$ws.Cells.Item(79, 9).Value = 0
$ws.Cells.Item(79, 10).Value = 1.58
$ws.Cells.Item(79, 11).Value = '28/12/2023 07:12'
$ws.Cells.Item(79, 12).Value = 1.39
$ws.Cells.Item(79, 13).Value = '29/12/2023 15:51'
$ws.Cells.Item(79, 14).Value = 3.32
$ws.Cells.Item(79, 15).Value = '28/12/2023 07:12'
$ws.Cells.Item(79, 16).Value = 4.13
$ws.Cells.Item(79, 17).Value = '29/12/2023 15:55'
$ws.Cells.Item(79, 18).Value = 6.24
$ws.Cells.Item(79, 19).Value = '28/12/2023 07:12'
$ws.Cells.Item(79, 20).Value = 10.82
$ws.Cells.Item(79, 21).Value = '29/12/2023 15:55'
$ws.Cells.Item(79, 22).Value = 'https://www.betexplorer.com/football/algeria/ligue-1/es-setif-oran/pI9HNfln/'

# Row 80
$ws.Range("A74:V74").Copy()
$ws.Range("A80:V80").PasteSpecial(-4122)
$ws.Cells.Item(80, 1).Value = 79
$ws.Cells.Item(80, 2).Value = 'algeria'
$ws.Cells.Item(80, 3).Value = 'ligue-1'
$ws.Cells.Item(80, 4).Value = '2023-2024'
$ws.Cells.Item(80, 5).Value = 45289.75
$ws.Cells.Item(80, 6).Value = 'USM Alger'
$ws.Cells.Item(80, 7).Value = 0
$ws.Cells.Item(80, 8).Value = 'MC Alger'
$ws.Cells.Item(80, 9).Value = 0
$ws.Cells.Item(80, 10).Value = 2.92
$ws.Cells.Item(80, 11).Value = '28/12/2023 07:12'
$ws.Cells.Item(80, 12).Value = 3.62
$ws.Cells.Item(80, 13).Value = '29/12/2023 17:59'
$ws.Cells.Item(80, 14).Value = 2.96
$ws.Cells.Item(80, 15).Value = '28/12/2023 07:12'
$ws.Cells.Item(80, 16).Value = 2.98
$ws.Cells.Item(80, 17).Value = '29/12/2023 17:56'
$ws.Cells.Item(80, 18).Value = 2.42
$ws.Cells.Item(80, 19).Value = '28/12/2023 07:12'
$ws.Cells.Item(80, 20).Value = 2.24
$ws.Cells.Item(80, 21).Value = '29/12/2023 17:59'
$ws.Cells.Item(80, 22).Value = 'https://www.betexplorer.com/football/algeria/ligue-1/usm-alger-mc-alger/jDDqjB3f/'

# Row 81
$ws.Range("A74:V74").Copy()
$ws.Range("A81:V81").PasteSpecial(-4122)
$ws.Cells.Item(81, 1).Value = 80
$ws.Cells.Item(81, 2).Value = 'algeria'
$ws.Cells.Item(81, 3).Value = 'ligue-1'
$ws.Cells.Item(81, 4).Value = '2023-2024'
$ws.Cells.Item(81, 5).Value = 45289.75
$ws.Cells.Item(81, 6).Value = 'Biskra'
$ws.Cells.Item(81, 7).Value = 1
$ws.Cells.Item(81, 8).Value = 'Ben Aknoun'
$ws.Cells.Item(81, 9).Value = 1
$ws.Cells.Item(81, 10).Value = 1.57
$ws.Cells.Item(81, 11).Value = '28/12/2023 07:12'
$ws.Cells.Item(81, 12).Value = 1.41
$ws.Cells.Item(81, 13).Value = '29/12/2023 17:57'
$ws.Cells.Item(81, 14).Value = 3.47
$ws.Cells.Item(81, 15).Value = '28/12/2023 07:12'
$ws.Cells.Item(81, 16).Value = 4.05
$ws.Cells.Item(81, 17).Value = '29/12/2023 17:57'
$ws.Cells.Item(81, 18).Value = 5.97
$ws.Cells.Item(81, 19).Value = '28/12/2023 07:12'
$ws.Cells.Item(81, 20).Value = 9.98
$ws.Cells.Item(81, 21).Value = '29/12/2023 17:57'
$ws.Cells.Item(81, 22).Value = 'https://www.betexplorer.com/football/algeria/ligue-1/biskra-es-ben-aknoun/z72UKhY4/'
